$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 (05-04-2015, Monday): add Time In/Out, change "late" hours from 1 to 0.5, change style to "late" (like row 8's original orange style) ---
# Do this first, before row 8's own style changes below.
$ws.Range("C14").Value = "07:04:16"
$ws.Range("D14").Value = "11:15:23"
$ws.Range("I14").Value = 0.5
$ws.Range("A8:P8").Copy()
$ws.Range("A14:P14").PasteSpecial(-4122)

# --- Row 7 (04-27-2015, Monday): add Time In/Out, clear "late" count, change style to plain (like row 12) ---
$ws.Range("C7").Value = "07:18:53"
$ws.Range("D7").Value = "19:25:46"
$ws.Range("I7").Value = ""
$ws.Range("A12:P12").Copy()
$ws.Range("A7:P7").PasteSpecial(-4122)

# --- Row 8 (04-28-2015, Tuesday): fix Time In, add Time Out, clear "late" count, change style to plain ---
$ws.Range("C8").Value = "07:08:19"
$ws.Range("D8").Value = "18:42:35"
$ws.Range("I8").Value = ""
$ws.Range("A12:P12").Copy()
$ws.Range("A8:P8").PasteSpecial(-4122)

# --- Row 9 (04-29-2015, Wednesday): add Time In/Out, clear "late" count, change style to plain ---
$ws.Range("C9").Value = "07:24:56"
$ws.Range("D9").Value = "18:59:59"
$ws.Range("I9").Value = ""
$ws.Range("A12:P12").Copy()
$ws.Range("A9:P9").PasteSpecial(-4122)

# --- Row 10 (04-30-2015, Thursday): add Time In/Out, add overtime hours, clear "late" count, change style to plain ---
$ws.Range("C10").Value = "07:07:59"
$ws.Range("D10").Value = "17:23:13"
$ws.Range("F10").Value = 1.25
$ws.Range("I10").Value = ""
$ws.Range("A12:P12").Copy()
$ws.Range("A10:P10").PasteSpecial(-4122)

$excel.CutCopyMode = 0
